$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "minhkhoi" lost all his money (bug/cheat fix) and gets a note in column I
$ws.Range("B4").Value = 0
$ws.Range("I4").Value = ";3;1;1;4"

# New player row appended to the table
$ws.Range("A6").Value = "superprovip"
$ws.Range("B6").Value = 750
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = ";0"
$ws.Range("G6").Value = ";32"
$ws.Range("H6").Value = ";-750.0"
